# Reorders the comma-separated "Recorded By" values in column G.
# The transformation applied to each cell's list of names/emails is:
#   1. Reverse the order of the items
#   2. Stable-sort the reversed items case-insensitively (so items that are
#      equal except for case keep their reversed relative order, while
#      differently-named items end up in case-insensitive ascending order)
# Cells with a single value (no comma) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val.GetType().Name -ne "String") { continue }
    if ($val.IndexOf(",") -lt 0) { continue }

    $parts = $val -split ","
    $trimmed = @()
    foreach ($p in $parts) {
        $trimmed += $p.Trim()
    }

    # Step 1: reverse the order of the items
    $reversed = @()
    for ($i = $trimmed.Count - 1; $i -ge 0; $i--) {
        $reversed += $trimmed[$i]
    }

    # Step 2: stable sort (case-insensitive) -- Sort-Object preserves the
    # relative order of items whose sort key compares equal.
    $indexed = @()
    foreach ($it in $reversed) {
        $indexed += @{ Key = $it.ToLower(); Val = $it }
    }
    $sorted = $indexed | Sort-Object -Property Key

    $newParts = @()
    foreach ($s in $sorted) {
        $newParts += $s.Val
    }

    $newVal = [string]::Join(", ", $newParts)

    if ($newVal -ne $val) {
        $cell.Value2 = $newVal
    }
}
